# Insert a new weekly record at row 39 ("Fruta / hortaliza, semanal").
# This pushes the existing rows 39..84 down to 40..85 and fills the
# newly opened row 39 with the new week's data. Columns that are constant
# across the whole data block (A,B,C,E,F,G,H,I,J,K,Q,T) are copied from the
# row directly above (row 38), which is untouched by the insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 39; existing data shifts down one row.
$ws.Rows("39:39").Insert()

# Columns that stay constant for every record in this block.
# (Read with .Value2 -- this COM shim's .Value getter doesn't resolve to the
#  underlying scalar; .Value2 does. The setter is fine with either.)
$ws.Range("A39").Value = $ws.Range("A38").Value2
$ws.Range("B39").Value = $ws.Range("B38").Value2
$ws.Range("C39").Value = $ws.Range("C38").Value2
$ws.Range("E39").Value = $ws.Range("E38").Value2
$ws.Range("F39").Value = $ws.Range("F38").Value2
$ws.Range("G39").Value = $ws.Range("G38").Value2
$ws.Range("H39").Value = $ws.Range("H38").Value2
$ws.Range("I39").Value = $ws.Range("I38").Value2
$ws.Range("J39").Value = $ws.Range("J38").Value2
$ws.Range("K39").Value = $ws.Range("K38").Value2
$ws.Range("Q39").Value = $ws.Range("Q38").Value2
$ws.Range("T39").Value = $ws.Range("T38").Value2

# New record-specific values.
$ws.Range("D39").Value = 44894
$ws.Range("L39").Value = "Primera"
$ws.Range("M39").Value = 220
$ws.Range("N39").Value = 3600
$ws.Range("O39").Value = 4000
$ws.Range("P39").Value = 3782
$ws.Range("R39").Value = "Provincia de Curicó"
$ws.Range("S39").Value = 1891
